$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1402.8
$ws.Range("I58").Value = 503.5
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 1510.5
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -1360.5
$ws.Range("N58").Value = -15300
$ws.Range("H107").Value = 2550.28
$ws.Range("I107").Value = 2445.5264
$ws.Range("J107").Value = 2882
$ws.Range("K107").Value = 2445.5264
$ws.Range("L107").Value = 2882
$ws.Range("M107").Value = -525.5264000000002
$ws.Range("N107").Value = -6722
$ws.Range("H113").Value = 3351.25
$ws.Range("I113").Value = 3502.5
$ws.Range("J113").Value = 3200
$ws.Range("K113").Value = 3502.5
$ws.Range("L113").Value = 3200
$ws.Range("M113").Value = -248.5
$ws.Range("N113").Value = -9708
$ws.Range("H132").Value = 2813.1765
$ws.Range("I132").Value = 3001.6667
$ws.Range("J132").Value = 1399.5
$ws.Range("K132").Value = 9005.000100000001
$ws.Range("L132").Value = 4198.5
$ws.Range("M132").Value = -6475.000100000001
$ws.Range("N132").Value = -9258.5
$ws.Range("H135").Value = 5363.7837
$ws.Range("I135").Value = 4280.52
$ws.Range("J135").Value = 7620.5835
$ws.Range("K135").Value = 38524.68000000001
$ws.Range("L135").Value = 68585.2515
$ws.Range("M135").Value = -35989.68000000001
$ws.Range("N135").Value = -73655.2515
$ws.Range("H137").Value = 7413458
$ws.Range("I137").Value = 11112270
$ws.Range("J137").Value = 15832.889
$ws.Range("K137").Value = 33336810
$ws.Range("L137").Value = 47498.667
$ws.Range("M137").Value = -33334260
$ws.Range("N137").Value = -52598.667

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 552221.75
$ws.Range("I32").Value = 621380
$ws.Range("J32").Value = 16245.25
$ws.Range("K32").Value = 621380
$ws.Range("L32").Value = 16245.25
$ws.Range("M32").Value = -621093
$ws.Range("N32").Value = -16819.25
$ws.Range("H63").Value = 2353.2
$ws.Range("I63").Value = 2353.2
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2353.2
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1667.2
$ws.Range("H66").Value = 2353.2
$ws.Range("I66").Value = 2353.2
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 11766
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -8334
$ws.Range("H74").Value = 1561478
$ws.Range("I74").Value = 2074895.8
$ws.Range("J74").Value = 21224.857
$ws.Range("K74").Value = 2074895.8
$ws.Range("L74").Value = 21224.857
$ws.Range("M74").Value = -2074021.8
$ws.Range("N74").Value = -22972.857
$ws.Range("H77").Value = 1561478
$ws.Range("I77").Value = 2074895.8
$ws.Range("J77").Value = 21224.857
$ws.Range("K77").Value = 10374479
$ws.Range("L77").Value = 106124.285
$ws.Range("M77").Value = -10370111
$ws.Range("N77").Value = -114860.285
$ws.Range("H80").Value = 41000
$ws.Range("I80").Value = 40000
$ws.Range("J80").Value = 42000
$ws.Range("K80").Value = 40000
$ws.Range("L80").Value = 42000
$ws.Range("M80").Value = -39002
$ws.Range("N80").Value = -43996
$ws.Range("H83").Value = 41000
$ws.Range("I83").Value = 40000
$ws.Range("J83").Value = 42000
$ws.Range("K83").Value = 120000
$ws.Range("L83").Value = 126000
$ws.Range("M83").Value = -115008
$ws.Range("N83").Value = -135984
$ws.Range("H132").Value = 6565.22
$ws.Range("I132").Value = 5110.778
$ws.Range("J132").Value = 7383.3438
$ws.Range("K132").Value = 15332.334
$ws.Range("L132").Value = 22150.0314
$ws.Range("M132").Value = -12802.334
$ws.Range("N132").Value = -27210.0314

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19797.182
$ws.Range("I82").Value = 13274.333
$ws.Range("J82").Value = 49150
$ws.Range("K82").Value = 13274.333
$ws.Range("L82").Value = 49150
$ws.Range("M82").Value = -12891.333
$ws.Range("N82").Value = -49916
$ws.Range("H85").Value = 19797.182
$ws.Range("I85").Value = 13274.333
$ws.Range("J85").Value = 49150
$ws.Range("K85").Value = 13274.333
$ws.Range("L85").Value = 49150
$ws.Range("M85").Value = -11948.333
$ws.Range("N85").Value = -51802
$ws.Range("H86").Value = 6743.878
$ws.Range("I86").Value = 6521.7666
$ws.Range("J86").Value = 7349.636
$ws.Range("K86").Value = 6521.7666
$ws.Range("L86").Value = 7349.636
$ws.Range("M86").Value = -5398.7666
$ws.Range("N86").Value = -9595.636
$ws.Range("H89").Value = 6743.878
$ws.Range("I89").Value = 6521.7666
$ws.Range("J89").Value = 7349.636
$ws.Range("K89").Value = 32608.833
$ws.Range("L89").Value = 36748.18
$ws.Range("M89").Value = -26992.833
$ws.Range("N89").Value = -47980.18
$ws.Range("H105").Value = 7554.5
$ws.Range("I105").Value = 6531.5557
$ws.Range("J105").Value = 16761
$ws.Range("K105").Value = 6531.5557
$ws.Range("L105").Value = 16761
$ws.Range("M105").Value = -4784.5557
$ws.Range("N105").Value = -20255

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 14519.23
$ws.Range("I105").Value = 12875
$ws.Range("J105").Value = 20000
$ws.Range("K105").Value = 12875
$ws.Range("L105").Value = 20000
$ws.Range("M105").Value = -11128
$ws.Range("N105").Value = -23494
$ws.Range("H107").Value = 887.13043
$ws.Range("I107").Value = 799.3333
$ws.Range("J107").Value = 1203.2
$ws.Range("K107").Value = 799.3333
$ws.Range("L107").Value = 1203.2
$ws.Range("M107").Value = 1120.6667
$ws.Range("N107").Value = -5043.2
$ws.Range("H122").Value = 13747.667
$ws.Range("I122").Value = 2838.1667
$ws.Range("J122").Value = 35566.668
$ws.Range("K122").Value = 8514.500100000001
$ws.Range("L122").Value = 106700.004
$ws.Range("M122").Value = -6064.500100000001
$ws.Range("N122").Value = -111600.004
$ws.Range("H132").Value = 4971
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 5132.8335
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 15398.5005
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -20458.5005
$ws.Range("H134").Value = 5112.175
$ws.Range("I134").Value = 3298.8
$ws.Range("J134").Value = 5716.6333
$ws.Range("K134").Value = 9896.400000000001
$ws.Range("L134").Value = 17149.8999
$ws.Range("M134").Value = -7361.400000000001
$ws.Range("N134").Value = -22219.8999

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 14289.714
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 14289.714
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 42869.142
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -44491.142
$ws.Range("H107").Value = 3740.56
$ws.Range("I107").Value = 1224.3636
$ws.Range("J107").Value = 5717.5713
$ws.Range("K107").Value = 3673.0908
$ws.Range("L107").Value = 17152.7139
$ws.Range("M107").Value = -1753.0908
$ws.Range("N107").Value = -20992.7139
$ws.Range("H117").Value = 1429
$ws.Range("I117").Value = 1009.6
$ws.Range("J117").Value = 1953.25
$ws.Range("K117").Value = 3028.8
$ws.Range("L117").Value = 5859.75
$ws.Range("M117").Value = 413.1999999999998
$ws.Range("N117").Value = -12743.75
$ws.Range("H130").Value = 9594.454
$ws.Range("I130").Value = 3107.8

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H126").Value = 2384.6316
$ws.Range("I126").Value = 2146.8667
$ws.Range("J126").Value = 3276.25
$ws.Range("K126").Value = 6440.6001
$ws.Range("L126").Value = 9828.75
$ws.Range("M126").Value = -3970.6001
$ws.Range("N126").Value = -14768.75
$ws.Range("H132").Value = 14119.03
$ws.Range("I132").Value = 14161.571
$ws.Range("J132").Value = 14087.685
$ws.Range("K132").Value = 42484.713
$ws.Range("L132").Value = 42263.055
$ws.Range("M132").Value = -39954.713
$ws.Range("N132").Value = -47323.055

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2084
$ws.Range("I16").Value = 2125
$ws.Range("J16").Value = 2002
$ws.Range("K16").Value = 2125
$ws.Range("L16").Value = 2002
$ws.Range("M16").Value = -1955
$ws.Range("N16").Value = -2342
$ws.Range("H22").Value = 2911.4546
$ws.Range("I22").Value = 2440.4285
$ws.Range("J22").Value = 3131.2666
$ws.Range("K22").Value = 2440.4285
$ws.Range("L22").Value = 3131.2666
$ws.Range("M22").Value = -2145.4285
$ws.Range("N22").Value = -3721.2666
$ws.Range("H27").Value = 2911.4546
$ws.Range("I27").Value = 2440.4285
$ws.Range("J27").Value = 3131.2666
$ws.Range("K27").Value = 2440.4285
$ws.Range("L27").Value = 3131.2666
$ws.Range("M27").Value = -2333.4285
$ws.Range("N27").Value = -3345.2666
$ws.Range("H40").Value = 5747.8
$ws.Range("I40").Value = 5747.8
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5747.8
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -5611.8
$ws.Range("H64").Value = 50150
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 50150
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 50150
$ws.Range("N64").Value = -50600
$ws.Range("H67").Value = 50150
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 50150
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 50150
$ws.Range("N67").Value = -51710
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1581.8889
$ws.Range("I113").Value = 1212.6
$ws.Range("J113").Value = 2043.5
$ws.Range("K113").Value = 3637.8
$ws.Range("L113").Value = 6130.5
$ws.Range("M113").Value = -1467.8
$ws.Range("N113").Value = -10470.5
$ws.Range("H132").Value = 55557556
$ws.Range("I132").Value = 83334340
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 250003020
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -250000490
$ws.Range("N132").Value = -17060
